# Update cryptocurrency price/volume data on sheet "cryptos"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.426.49"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.671.72"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "645.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.497"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.442"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000233"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "4.290.98"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "3.665.52"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "69.404.13"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").Value = "3.817.53"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.46"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("D36").Value = "3.661.29"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.91"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.01"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000271"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "27.12"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.69%  "
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.42%  "
